# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data block (rows 12-13),
# pushing the existing rows 12:35 down to 14:37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 12.
$ws.Rows("12:13").Insert()

# Row 12 - new record (Primera)
$ws.Cells.Item(12, 1).Value2 = 5
$ws.Cells.Item(12, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(12, 3).Value2 = "Maule"
$ws.Cells.Item(12, 4).Value2 = 44483
$ws.Cells.Item(12, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12, 5).Value2 = 7
$ws.Cells.Item(12, 6).Value2 = "Fruta"
$ws.Cells.Item(12, 7).Value2 = 100107
$ws.Cells.Item(12, 8).Value2 = "Otros"
$ws.Cells.Item(12, 9).Value2 = 100107002
$ws.Cells.Item(12, 10).Value2 = "Chirimoya"
$ws.Cells.Item(12, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(12, 12).Value2 = "Primera"
$ws.Cells.Item(12, 13).Value2 = 80
$ws.Cells.Item(12, 14).Value2 = 25000
$ws.Cells.Item(12, 15).Value2 = 25000
$ws.Cells.Item(12, 16).Value2 = 25000
$ws.Cells.Item(12, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(12, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(12, 19).Value2 = 2500
$ws.Cells.Item(12, 20).Value2 = 10

# Row 13 - new record (Segunda)
$ws.Cells.Item(13, 1).Value2 = 5
$ws.Cells.Item(13, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(13, 3).Value2 = "Maule"
$ws.Cells.Item(13, 4).Value2 = 44483
$ws.Cells.Item(13, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13, 5).Value2 = 7
$ws.Cells.Item(13, 6).Value2 = "Fruta"
$ws.Cells.Item(13, 7).Value2 = 100107
$ws.Cells.Item(13, 8).Value2 = "Otros"
$ws.Cells.Item(13, 9).Value2 = 100107002
$ws.Cells.Item(13, 10).Value2 = "Chirimoya"
$ws.Cells.Item(13, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(13, 12).Value2 = "Segunda"
$ws.Cells.Item(13, 13).Value2 = 30
$ws.Cells.Item(13, 14).Value2 = 22000
$ws.Cells.Item(13, 15).Value2 = 22000
$ws.Cells.Item(13, 16).Value2 = 22000
$ws.Cells.Item(13, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(13, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(13, 19).Value2 = 2200
$ws.Cells.Item(13, 20).Value2 = 10
